# tree2l.xlsx fix: shift the per-internode "MTG code" label cells two
# columns to the left (D -> B, E -> C) for the data block in rows 40-122,
# and refresh the saved sheet view (scroll position / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A2B")

# Rows whose column-D cell must move to column B.
$dRows = @(40,41,44,47,50,53,56,59,62,65,68,71,74,77,80,83,86,89,92,95,98,101,104,107,110,113,116,119,122)

# Rows whose column-E cell must move to column C.
$eRows = @(42,43,45,46,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,64,65,66,67,68,69,70,72,73,75,76,78,79,81,82,84,85,87,88,90,91,93,94,96,97,99,100,102,103,105,106,108,109,111,112,114,115,117,118,120,121)

foreach ($r in $dRows) {
    $src = $ws.Cells.Item($r, 4)   # column D
    $dst = $ws.Cells.Item($r, 2)   # column B
    $src.Cut($dst) | Out-Null
    $src.Clear() | Out-Null
}

foreach ($r in $eRows) {
    $src = $ws.Cells.Item($r, 5)   # column E
    $dst = $ws.Cells.Item($r, 3)   # column C
    $src.Cut($dst) | Out-Null
    $src.Clear() | Out-Null
}

# Restore the sheet view: scrolled to A11, with E40 selected.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E40").Select() | Out-Null
